$excel.DisplayAlerts = $false

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)      # invalidCredentialTest

# ------------------------------------------------------------------
# "validCredentialTest" used to be blank. In the real edit it was
# produced by duplicating the invalidCredentialTest tab (hence it
# shares the exact same best-fit column widths) and then swapping in
# new test data. Reproduce that: drop the blank sheet, copy sheet1
# into its place, rename it, and overwrite the cell values.
# ------------------------------------------------------------------
$wb.Worksheets.Item(2).Delete() | Out-Null
$ws1.Copy($null, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "validCredentialTest"

$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"
$ws2.Range("C1").Value = "Language"
$ws2.Range("D1").Value = "Expected Value"

$ws2.Range("A2").Value = "admin"
$ws2.Range("B2").Value = "pass"
$ws2.Range("C2").Value = "English (Indian)"
$ws2.Range("D2").Value = "OpenEMR"

$ws2.Range("A3").Value = "accountant"
$ws2.Range("B3").Value = "accountant"
$ws2.Range("C3").Value = "Dutch"
$ws2.Range("D3").Value = "OpenEMR"

# invalidCredentialTest: plain A1:D3 selection, no longer the visible tab
$ws1.Select() | Out-Null
$ws1.Range("A1:D3").Select() | Out-Null

# validCredentialTest becomes the active/visible tab, selection D2:D3
$ws2.Select() | Out-Null
$ws2.Range("D2:D3").Select() | Out-Null
